# Principal Component Analysis.docx edits
#
# 1) Merge the three runs around "e.g." (which were split apart by
#    <w:proofErr w:type="gramStart"/>...<w:proofErr w:type="gramEnd"/>) into a
#    single run reading "...reduce dimensions of the data e.g. " by doing a
#    Find/Replace of the phrase with itself - Word's Find engine rewrites the
#    matched span as one run and drops the (now stale) grammar-check markup.
# 2) Append a new run containing "#" right after "...has maximal variance ".
# 3) Add two empty paragraphs and a new paragraph asking "Why do we choose
#    the particular constaint of sum of all coefficients squared being 1".

$d = $word.ActiveDocument

# --- Step 1: clean up the "e.g." runs / drop the gramStart/gramEnd proofErr markup ---
$d.Content.Find.Execute(
    "reduce dimensions of the data e.g. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "reduce dimensions of the data e.g. ", 2) | Out-Null

# --- Step 2: append a "#" run after "...has maximal variance " ---
$findRange = $d.Content
$findRange.Find.Execute(
    " new variable which is easier to work with and has maximal variance ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$variancePara = $d.Range($findRange.Start, $findRange.End)
$variancePkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p>' + `
  '<w:r><w:t xml:space="preserve"> new variable which is easier to work with and has maximal variance </w:t></w:r>' + `
  '<w:r><w:t>#</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$variancePara.InsertXML($variancePkg) | Out-Null

# --- Step 3: append two empty paragraphs and the new "constaint" paragraph ---
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$wholeLastPara = $d.Range($lastRange.Start, $lastRange.End)

$tailPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p>' + `
  '<w:r><w:t xml:space="preserve">Can use this weighted combination instead of the actual variables to reduce dimensions of the data e.g. </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>intoruce</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> new variable which is easier to work with and has maximal variance </w:t></w:r>' + `
  '<w:r><w:t>#</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p/>' + `
  '<w:p/>' + `
  '<w:p>' + `
  '<w:r><w:t xml:space="preserve">Why do we choose the particular </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>constaint</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> of sum of all coefficients squared being 1  </w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$wholeLastPara.InsertXML($tailPkg) | Out-Null

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
